$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2:C19").Value = 8347
$ws.Range("C20:C31").Value = 7752
$ws.Range("C32:C34").Value = 7694
$ws.Range("C35:C50").Value = 7610
$ws.Range("C51:C252").Value = 7534
